# Data.xlsx / "Base" sheet update
#
# 1) FINBP247 (Sunil K, row 4) has finished being rated -> Category flips
#    from "TBD" to "New to Rate".
# 2) The block of rows 49-79 (the old "Rahul Singh" / "Foram Jigar Buch"
#    reportees that had been pasted in twice / were stale) is cleared out,
#    shrinking the used range from A1:Q79 down to A1:Q77.
# 3) The sheet view scrolls back up and the active cell moves to D48.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Category correction for row 4 --------------------------------
$ws.Range("H4").Value = "New to Rate"

# --- 2) Wipe out the old trailing rows (49 through 79) ----------------
$ws.Rows("49:79").ClearContents()

# --- 3) Restore the view: scroll position + active cell ---------------
$ws.Activate()
$ws.Range("D48").Select()
